$d = $word.ActiveDocument

# --- Merge the split "Text" + "N" runs into single "TextN" runs ---
# These were originally two separate <w:r> runs in the same paragraph
# ("Text" and "2"/"3"/"4"); a Find & Replace across the run boundary
# collapses them into a single run while keeping identical formatting.
$d.Content.Find.Execute("Text2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Text2", 2) | Out-Null
$d.Content.Find.Execute("Text3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Text3", 2) | Out-Null
$d.Content.Find.Execute("Text4", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Text4", 2) | Out-Null

# --- Append a new row to the Data Dictionary table for "Clock" ---
$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()
$newRow.Height = 13.05

$newRow.Cells.Item(1).Range.Text = "Clock"
$newRow.Cells.Item(2).Range.Text = "integer"
$newRow.Cells.Item(3).Range.Text = "Pygame.time.clock()"
$newRow.Cells.Item(4).Range.Text = "This is the pygame clock class use to control the frame rate "
